$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 (PSG 18) values
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 4

# Update selection to column E (entire column selected)
$ws.Columns("E").Select()
